# Add the new "Class Information" worksheet as the last sheet in the workbook.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Class Information"

# Header row
$ws.Cells.Item(1, 1).Value = "Class Name"
$ws.Cells.Item(1, 2).Value = "Time Spend"
$ws.Cells.Item(1, 3).Value = "Class ID"
$ws.Cells.Item(1, 4).Value = "Professor"
$ws.Cells.Item(1, 5).Value = "Room"
$ws.Cells.Item(1, 6).Value = "Max Occupancy"
$ws.Cells.Item(1, 7).Value = "Current Occupancy"

# Class data rows: Class Name, Time Spend, Class ID, Professor, Room, Max Occupancy
$rows = @(
    @("English",        "8:00-9:00",   12345, "Jane Doe",        "A-123", 30),
    @("Geometry",        "9:00-10:00",  12346, "Kevin Smith",     "A-231", 30),
    @("Physics",         "10:00-11:00", 12347, "Jennifer Jones",  "B-102", 30),
    @("Chemistry",       "11:00-12:00", 12348, "David Smith",     "B-204", 25),
    @("Biology",         "12:00-13:00", 12349, "Eric Huang",      "B-123", 25),
    @("Philosphy",       "13:00-14:00", 12350, "Vincent Ku",      "A-123", 30),
    @("Linear Algebra",  "14:00-15:00", 12351, "Joshua Hu",       "A-231", 30),
    @("Discrete Math",   "15:00-16:00", 12352, "Linda J",         "B-102", 30),
    @("Calculus",        "16:00-17:00", 12353, "Henry H",         "B-204", 25),
    @("Python",          "17:00-18:00", 12354, "Julie J",         "B-123", 25),
    @("Java",            "18:00-19:00", 12355, "Jason L",         "A-210", 30)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Row 5 ("Chemistry") has its Time Spend cell formatted as a time value.
$ws.Cells.Item(5, 2).NumberFormat = "h:mm"

# Make the new sheet the active tab, with the same selection state captured in the diff.
$ws.Activate()
$ws.Range("H21").Select() | Out-Null
